# DataPool.xlsx edit script
# Commit: "Se agrega Metodo EditarInstitucion en clase PageInstituciones //
#          Se agrega script 0124 a la clase Tests_AdmInstituciones"
#
# Effect on the "Hoja1" worksheet data table:
#   - Row 41 (DEC_0124) keeps its institution code but the rest of the row
#     (D:J -- prefix/description/email/name/etc.) is cleared, since the
#     "Nueva Empresa QA 124" scenario row no longer carries that data.
#   - Row 42 changes from DEC_0125 / "Nueva Empresa QA 125" down to just
#     DEC_0125 with the rest of the row cleared (same treatment as row 41).
#   - Six new script rows are appended (DEC_0126 .. DEC_0131), each only
#     populated in columns A (code), B and C (fixed constants), mirroring
#     the now-empty placeholder rows that used to follow.
#   - The previously-blank row that used to sit at row 53 is removed so
#     that everything below shifts up by one row (used range goes from
#     A1:J55 to A1:J54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank spacer row - this shifts rows 54/55 up to 53/54 and
# updates the sheet dimension/used range accordingly.
$ws.Rows("53:53").Delete()

# Row 41: keep the DEC_0124 code in column A, clear everything else.
$ws.Range("D41:G41").Clear()
$ws.Range("H41").ClearContents()
$ws.Range("I41:J41").Clear()

# Row 42: becomes DEC_0125 (previously held DEC_0125 data one row lower);
# clear the rest of the row just like row 41.
$ws.Range("A42").Value = "DEC_0125"
$ws.Range("D42:G42").Clear()
$ws.Range("H42").ClearContents()
$ws.Range("I42:J42").Clear()

# New script rows DEC_0126 .. DEC_0131 in rows 43-48, each only filling in
# columns A (code), B and C (constant TC credentials); column H keeps the
# same blank/formatted placeholder cell the row already had.
$newCodes = @("DEC_0126", "DEC_0127", "DEC_0128", "DEC_0129", "DEC_0130", "DEC_0131")
$row = 43
foreach ($code in $newCodes) {
    $ws.Range("A$row").Value = $code
    $ws.Range("B$row").Value = "13712759-8"
    $ws.Range("C$row").Value = "Verity1.0"
    $row = $row + 1
}
